# Mark additional test script rows as "Execute = Yes" in the customer test script workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rows 2-66, 68-76, 79-89 and 97-102 flip their "Execute" flag (column B) from "No" to "Yes".
# Rows 67, 77, 78 and 90-96 are already "Yes" and stay unchanged.
$ws.Range("B2:B66").Value = "Yes"
$ws.Range("B68:B76").Value = "Yes"
$ws.Range("B79:B89").Value = "Yes"
$ws.Range("B97:B102").Value = "Yes"
